# Logboek update: hours for rows 10-12 were reshuffled (day 10 is now the
# empty/absent day, day 11 now carries the hours that used to sit on day 10)
# and day 12's end-time / break were corrected slightly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 10 -> now empty (no hours logged)
$ws.Range("B10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("F10").Value = 0

# Row 11 -> now carries the hours that used to be on row 10
$ws.Range("B11").Value = 9
$ws.Range("D11").Value = 14
$ws.Range("F11").Value = 1

# Row 12 -> end time and break corrected
$ws.Range("D12").Value = 14
$ws.Range("F12").Value = 0

# Update the window view/selection to match where the user left off editing
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("I12:I13").Select()
